$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "salted" output values (column B), replacing the prior random-salt run.
$newValues = @(4156,-1353,2501,1368,-4059,3499,-1738,2909,3706,-898,-2041,-1323,-2314,3306,-1508,6676,829,4632,-599,5709,-382,-300,42,2564,81,2014,2937,8623,9968,10486,4384,8760,5912,5577,5257,7717,14310,14186,5704,9438,15610,16517,14283,10850,9791,9851,12543,18889,12299,14981,15156,21243,16933,16720,25166,17304,18973,21304,20052,28739,27228,31262,28500,24500,33637,32072,28999,33161,36752,35689,38143,31458,38633,41821,40771,41197,45377,40082,41162,45697,41396,44137,45337,51169,51172,49695,55770,55884,53568,59990,60625,54437,55344,59865,66895,66477,68535,68905,70196,72986)

$arr = New-Object 'object[,]' 100,1
for ($i = 0; $i -lt 100; $i++) {
    $arr[$i, 0] = $newValues[$i]
}
$ws.Range("B1:B100").Value = $arr

# Update the chart title text (Excel collapses the two runs it previously had
# into a single run when the title text is assigned as a whole string).
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Salted Output (Range: 1000 - 5000)"
